$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the B5 value (gender label) - was incorrectly "أنثى" (Female), should be "الذكر" (Male)
$ws.Range("B5").Value = "الذكر"

# Set column B width (matches bestFit width observed in target)
$ws.Columns.Item(2).ColumnWidth = 6.43

# Set the active cell / selection
$ws.Range("D16").Select()

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
